$d = $word.ActiveDocument

$find0 = "Siña tokante e famia, yu i relashon entre mayor i yu."
$replace0 = "Siña tokante e famia, yu i relashon entre mayor/edukadó i yu."
$found0 = $d.Content.Find.Execute($find0, $true, $true, $false, $false, $false, $true, 1, $false, $replace0, 2)
if (-not $found0) { Write-Host "NOT FOUND: 0" }

$find1 = "Bo mester usa konsultanan individual komo un oportunidat pa siña mas tantu ku bo por tokante e sirkunstansianan i retonan di e famia. Esaki ta importante ya ku ora un mayor kuminsá kambia e manera ku e ta interkambiá ku e yunan den e kas, e lo afektá tur otro hende i tin bia por forma konflikto. "
$replace1 = "Bo mester usa konsultanan individual komo un oportunidat pa siña mas tantu ku bo por tokante e sirkunstansianan i retonan di e famia. Esaki ta importante ya ku ora un mayor/edukadó kuminsá kambia e manera ku e ta interkambiá ku e yunan den e kas, e lo afektá tur otro hende i tin bia por forma konflikto. "
$found1 = $d.Content.Find.Execute($find1, $true, $true, $false, $false, $false, $true, 1, $false, $replace1, 2)
if (-not $found1) { Write-Host "NOT FOUND: 1" }

$find2 = "E kombersashon aki lo duná bo un komprondementu mas grandi di algun di e retonan ku famianan ta eksperensiá na kas. E lo ekipá bo tambe ku informashon ku por ta útil pa yuda mayornan establesé metanan semanal i solushoná asuntunan ku ta surgi ora mayornan ta praktiká abilidatnan nobo na kas."
$replace2 = "E kombersashon aki lo lagá bo komprondé mihó algun di e retonan ku famianan ta eksperensiá na kas. E lo ekipá bo tambe ku informashon ku por ta útil pa yuda mayornan/edukadónan establesé metanan semanal i solushoná asuntunan ku ta surgi ora mayornan/edukadónan ta praktiká abilidatnan nobo na kas."
$found2 = $d.Content.Find.Execute($find2, $true, $true, $false, $false, $false, $true, 1, $false, $replace2, 2)
if (-not $found2) { Write-Host "NOT FOUND: 2" }

$find3 = "Tambe ta importante pa bo yuda mayornan identifiká UN meta positivo, spesífiko i realístiko pa e yu ku nan ta enfoká riba dje durante e programa."
$replace3 = "Tambe ta importante pa bo yuda mayornan/edukadónan identifiká UN meta positivo, spesífiko i realístiko pa e yu ku nan ta enfoká riba dje durante e programa."
$found3 = $d.Content.Find.Execute($find3, $true, $true, $false, $false, $false, $true, 1, $false, $replace3, 2)
if (-not $found3) { Write-Host "NOT FOUND: 3" }

$find4 = "Pa yuda mayornan pone un meta, ta importante pa bo pidi mayornan pa deskribí nan ekspektativanan tokante e programa. Ora bo ta hasi esaki, yuda nan identifiká un meta spesífiko tokante kon nan lo ke kambia nan relashon ku nan yu òf e komportashon di nan yu. "
$replace4 = "Pa yuda mayornan/edukadónan pone un meta, ta importante pa bo pidi mayornan/edukadónan pa deskribí nan ekspektativanan tokante e programa. Ora bo ta hasi esaki, yuda nan identifiká un meta spesífiko tokante kon nan lo ke kambia nan relashon ku nan yu òf e komportashon di nan yu. "
$found4 = $d.Content.Find.Execute($find4, $true, $true, $false, $false, $false, $true, 1, $false, $replace4, 2)
if (-not $found4) { Write-Host "NOT FOUND: 4" }

$find5 = "Hopi bia mayornan lo duna metanan vago manera, “Mi ke pa mi yu ta bon mucha,” òf “Mi ke pa mi yu tin éksito na skol.” Kisas bo mester guia nan pa ta mas spesífiko tokante kiko e mayor ke men ku “bon mucha” òf “tin éksito na skol.” "
$replace5 = "Hopi bia mayornan/edukadónan lo duna metanan vago manera, “Mi ke pa mi yu ta bon mucha,” òf “Mi ke pa mi yu tin éksito na skol.” Kisas bo mester guia nan pa ta mas spesífiko tokante kiko e mayor/edukadó ke men ku “bon mucha” òf “tin éksito na skol.” "
$found5 = $d.Content.Find.Execute($find5, $true, $true, $false, $false, $false, $true, 1, $false, $replace5, 2)
if (-not $found5) { Write-Host "NOT FOUND: 5" }

$find6 = "Ta bo trabou pa yuda mayornan deskribí un komportashon mas spesífiko posibel. Lo bo mester usa pregunta pa logra ku nan ta bira mas i mas spesífiko. Nos ta rekomendá pa bo puntra mayornan kiko “ta bon mucha” òf “tin éksito na skol” ta nifiká pa nan. Ademas, eksplorando ki ora, unda i dikon nan lo ke pa e komportashon akí sosodé, lo yuda hasi nan meta mas spesífiko."
$replace6 = "Ta bo trabou pa yuda mayornan/edukadónan deskribí un komportashon mas spesífiko posibel. Lo bo mester usa pregunta pa logra ku nan ta bira mas i mas spesífiko. Nos ta rekomendá pa bo puntra mayornan/edukadónan kiko “ta bon mucha” òf “tin éksito na skol” ta nifiká pa nan. Ademas, si bo eksplorá ki ora, unda i dikon nan lo ke pa e komportashon akí sosodé, lo yuda hasi nan meta mas spesífiko."
$found6 = $d.Content.Find.Execute($find6, $true, $true, $false, $false, $false, $true, 1, $false, $replace6, 2)
if (-not $found6) { Write-Host "NOT FOUND: 6" }

$find7 = "Pone metanan spesífiko lo yuda e mayornan identifiká komportashonnan ku nan por yuda enkurashá pa sosodé mas tantu segun ku nan ta desaroyá relashonnan mas positivo ku nan yu. "
$replace7 = "Pone metanan spesífiko lo yuda e mayornan/edukadónan identifiká komportashonnan ku nan por yuda stimulá mas tantu segun ku nan ta desaroyá relashonnan mas positivo ku nan yu. "
$found7 = $d.Content.Find.Execute($find7, $true, $true, $false, $false, $false, $true, 1, $false, $replace7, 2)
if (-not $found7) { Write-Host "NOT FOUND: 7" }

$find8 = "Un mayor ku ke pa su yu “ta bon mucha” eventualmente por tin un meta manera “Mi ke pa mi yu kuminda su grandinan na un manera respetuoso ora e drenta kas.” "
$replace8 = "Un mayor/edukadó ku ke pa su yu “ta bon mucha” eventualmente por tin un meta manera “Mi ke pa mi yu kuminda su grandinan na un manera respetuoso ora e drenta kas.” "
$found8 = $d.Content.Find.Execute($find8, $true, $true, $false, $false, $false, $true, 1, $false, $replace8, 2)
if (-not $found8) { Write-Host "NOT FOUND: 8" }

$find9 = "Di mes manera, bo por yuda un mayor ku ke pa su yu “tin éksito na skol” pa identifiká un meta mas spesífiko manera, “Mi ke pa mi yu hasi su hùiswèrk mésora ku e yega kas for di skol.” "
$replace9 = "Di mes manera, bo por yuda un mayor/edukadó ku ke pa su yu “tin éksito na skol” pa identifiká un meta mas spesífiko manera, “Mi ke pa mi yu hasi su hùiswèrk mésora ku e yega kas for di skol.” "
$found9 = $d.Content.Find.Execute($find9, $true, $true, $false, $false, $false, $true, 1, $false, $replace9, 2)
if (-not $found9) { Write-Host "NOT FOUND: 9" }

$find10 = "Por ehèmpel, en bes di bisa, “Mi ke pa mi yu stòp di papia palabra mahos ku mi,” mester yuda un mayor pa menshoná e komportashon ku e ke mira: “Mi ke pa mi yu usa palabranan amabel ora e ta papia ku mi.” "
$replace10 = "Por ehèmpel, en bes di bisa, “Mi ke pa mi yu stòp di papia palabra mahos ku mi,” mester yuda un mayor/edukadó pa menshoná e komportashon ku e ke mira: “Mi ke pa mi yu usa palabranan amabel ora e ta papia ku mi.” "
$found10 = $d.Content.Find.Execute($find10, $true, $true, $false, $false, $false, $true, 1, $false, $replace10, 2)
if (-not $found10) { Write-Host "NOT FOUND: 10" }

$find11 = "Ora abo i e mayor ta kontentu ku e meta spesífiko"
$replace11 = "Ora abo i e mayor/edukadó ta kontentu ku e meta spesífiko"
$found11 = $d.Content.Find.Execute($find11, $true, $true, $false, $false, $false, $true, 1, $false, $replace11, 2)
if (-not $found11) { Write-Host "NOT FOUND: 11" }

$find12 = "Mayornan hopi bia lo ke skohe metanan ku no ta alkansabel durante e programa òf tin bia ta imposibel pa e fase di desaroyo di nan yu. Bo por guia mayornan ku pasenshi i komprenshon pa tin espektativanan mas realístiko pa nan yunan."
$replace12 = "Mayornan/Edukadónan hopi bia lo ke skohe metanan ku no ta alkansabel durante e programa òf tin bia ta imposibel pa e yu por logra segun su fase di desaroyo. Bo por guia mayornan/edukadónan ku pasenshi i komprenshon pa tin espektativanan mas realístiko pa nan yunan."
$found12 = $d.Content.Find.Execute($find12, $true, $true, $false, $false, $false, $true, 1, $false, $replace12, 2)
if (-not $found12) { Write-Host "NOT FOUND: 12" }

$find13 = "Por ehèmpel, si un mayor bisa ku e ke pa su yu pasa su èksamennan di skol na fin di aña, bo por puntr’é ki komportashonnan spesífiko su yu mester hasi pa hasi esei posibel."
$replace13 = "Por ehèmpel, si un mayor/edukadó bisa ku e ke pa su yu pasa su èksamennan di skol ku éksito na fin di aña, bo por puntr’é ki komportashonnan spesífiko su yu mester hasi pa hasi esei posibel."
$found13 = $d.Content.Find.Execute($find13, $true, $true, $false, $false, $false, $true, 1, $false, $replace13, 2)
if (-not $found13) { Write-Host "NOT FOUND: 13" }

$find14 = "Di mes manera, si un mayor ke pa su yu di 2 aña por bisti su so mainta, bo mester eksplorá si esei ta un ekspektativa rasonabel i dikon e yu no por bisti su so. Despues bo por yuda e mayor identifiká un meta mas realístiko ku lo yuda su yu desaroyá e abilidatnan pa bisti su so manera, “Mi lo ke pa mi yu koperá ku mi miéntras mi ta yud’é bisti paña mainta.” "
$replace14 = "Di mes manera, si un mayor/edukadó ke pa su yu di 2 aña por bisti su paña su so mainta, bo mester eksplorá si esei ta un ekspektativa rasonabel i dikon e yu no por bisti su so. Despues bo por yuda e mayor/edukadó identifiká un meta mas realístiko ku lo yuda su yu desaroyá e abilidatnan pa bisti paña su so manera, “Mi lo ke pa mi yu koperá ku mi miéntras mi ta yud’é bisti paña mainta.” "
$found14 = $d.Content.Find.Execute($find14, $true, $true, $false, $false, $false, $true, 1, $false, $replace14, 2)
if (-not $found14) { Write-Host "NOT FOUND: 14" }

$find15 = "Finalmente, konsultanan individual ta oportunidatnan pa papia tokante kualke asuntu logístiko ku e mayornan, tokante e ora pa e seshonnan di grupo, akseso na un telefòn selular, kualke nesesidat di rekargá data i preguntanan teknológiko."
$replace15 = "Finalmente, konsultanan individual ta oportunidatnan pa papia tokante kualke asuntu logístiko ku e mayornan/edukadónan, tokante e ora pa e seshonnan di grupo, akseso na un telefòn selular, kualke nesesidat di rekargá data i preguntanan teknológiko."
$found15 = $d.Content.Find.Execute($find15, $true, $true, $false, $false, $false, $true, 1, $false, $replace15, 2)
if (-not $found15) { Write-Host "NOT FOUND: 15" }

$find16 = "A. Introdusí boso mes na e mayor i henter famia si nan ta presente."
$replace16 = "A. Introdusí boso mes na e mayor/edukadó i henter famia si nan ta presente."
$found16 = $d.Content.Find.Execute($find16, $true, $true, $false, $false, $false, $true, 1, $false, $replace16, 2)
if (-not $found16) { Write-Host "NOT FOUND: 16" }

$find17 = "Ken mas ta biba einan? Kuantu mucha tin? Esposo/Esposa? Pareha? Grandinan?"
$replace17 = "Ken mas ta biba einan? Kuantu mucha tin? Esposo/Esposa? Pareha? Grandinan (Wela, Tawela)?"
$found17 = $d.Content.Find.Execute($find17, $true, $true, $false, $false, $false, $true, 1, $false, $replace17, 2)
if (-not $found17) { Write-Host "NOT FOUND: 17" }

$find18 = "Ken mas ta duna kuido na mucha?"
$replace18 = "Ken mas ta kuida e muchanan?"
$found18 = $d.Content.Find.Execute($find18, $true, $true, $false, $false, $false, $true, 1, $false, $replace18, 2)
if (-not $found18) { Write-Host "NOT FOUND: 18" }

$find19 = "D. Kombersá ku e mayor tokante su relashon ku su yu:"
$replace19 = "D. Kombersá ku e mayor/edukadó tokante su relashon ku su yu:"
$found19 = $d.Content.Find.Execute($find19, $true, $true, $false, $false, $false, $true, 1, $false, $replace19, 2)
if (-not $found19) { Write-Host "NOT FOUND: 19" }

$find20 = "Kòrda e mayor ku e lo enfoká riba e yu ku el a skohe durante di e"
$replace20 = "Kòrda e mayor/edukadó ku e lo enfoká riba e yu ku el a skohe durante di e"
$found20 = $d.Content.Find.Execute($find20, $true, $true, $false, $false, $false, $true, 1, $false, $replace20, 2)
if (-not $found20) { Write-Host "NOT FOUND: 20" }

$find21 = "Si e mayor no a selektá un yu spesífiko pa enfoká riba dje durante e programa: "
$replace21 = "Si e mayor/edukadó no a selektá un yu spesífiko pa enfoká riba dje durante e programa: "
$found21 = $d.Content.Find.Execute($find21, $true, $true, $false, $false, $false, $true, 1, $false, $replace21, 2)
if (-not $found21) { Write-Host "NOT FOUND: 21" }

$find22 = "Pidi e mayor pa selektá un yu pa enfoká riba dje durante e programa. E yu akí mester ta entre e edat di 2 i 17 aña. Si e mayor tin mas ku un yu entre e rango di edat akí, e mester selektá e yu ku aworaki e tin e relashon òf retonan mas difísil pa manehá su komportashon. "
$replace22 = "Pidi e mayor/edukadó pa selektá un yu pa enfoká riba dje durante e programa. E yu akí mester ta entre e edat di 2 pa 17 aña. Si e mayor/edukadó tin mas ku un yu entre e rango di edat akí, e mester selektá e yu ku aworaki e tin e relashon mas difísil kuné òf ku ta dun,é mas reto pa manehá su komportashon. "
$found22 = $d.Content.Find.Execute($find22, $true, $true, $false, $false, $false, $true, 1, $false, $replace22, 2)
if (-not $found22) { Write-Host "NOT FOUND: 22" }

$find23 = "Bo por sigurá e mayor tambe ku e por apliká e abilidatnan ku e ta siña den e programa pa tur e otro muchanan den su famia, pero ku e mester enfoká riba e mucha akí durante e diskushonnan di grupo i práktika na kas."
$replace23 = "Bo por sigurá e mayor/edukadó tambe ku e por apliká e abilidatnan ku e ta siña den e programa pa tur e otro muchanan den su famia, pero ku e mester enfoká riba e mucha akí durante e diskushonnan di grupo i práktika na kas."
$found23 = $d.Content.Find.Execute($find23, $true, $true, $false, $false, $false, $true, 1, $false, $replace23, 2)
if (-not $found23) { Write-Host "NOT FOUND: 23" }

$find24 = "Kua ta algun reto ku bo ta konfrontando pa manehá e komportashon di bo yu?"
$replace24 = "Kua ta algun reto ku bo ta konfrontá pa manehá e komportashon di bo yu?"
$found24 = $d.Content.Find.Execute($find24, $true, $true, $false, $false, $false, $true, 1, $false, $replace24, 2)
if (-not $found24) { Write-Host "NOT FOUND: 24" }

$find25 = "Tin otro retonan ku ta hasié difísil pa ta un mayor?"
$replace25 = "Tin otro retonan ku ta hasié difísil pa ta un mayor/edukadó?"
$found25 = $d.Content.Find.Execute($find25, $true, $true, $false, $false, $false, $true, 1, $false, $replace25, 2)
if (-not $found25) { Write-Host "NOT FOUND: 25" }

$find26 = "E. Metanan di Mayor pa e Programa"
$replace26 = "E. Metanan di Mayor/Edukadó pa e Programa"
$found26 = $d.Content.Find.Execute($find26, $true, $true, $false, $false, $false, $true, 1, $false, $replace26, 2)
if (-not $found26) { Write-Host "NOT FOUND: 26" }

$find27 = "Yuda e mayor identifiká UN meta spesífiko, positivo i realístiko."
$replace27 = "Yuda e mayor/edukadó identifiká UN meta spesífiko, positivo i realístiko."
$found27 = $d.Content.Find.Execute($find27, $true, $true, $false, $false, $false, $true, 1, $false, $replace27, 2)
if (-not $found27) { Write-Host "NOT FOUND: 27" }

$find28 = "Skibi e meta di e mayor na parti abou di bo profil di partisipante/mayor."
$replace28 = "Skibi e meta di e mayor/edukadó na parti abou di bo profil di partisipante/mayor/edukadó."
$found28 = $d.Content.Find.Execute($find28, $true, $true, $false, $false, $false, $true, 1, $false, $replace28, 2)
if (-not $found28) { Write-Host "NOT FOUND: 28" }

$find29 = "Aklará e nivel di alfabetisashon di mayornan – lo bo mester manda mensahenan di oudio en bes di mensahenan di teksto si mayornan tin difikultat pa lesa. Lo bo por mester splika kon un partisipante ta skucha un mensahe di oudio."
$replace29 = "Aklará e nivel di alfabetisashon di mayornan/edukadónan – lo bo mester manda mensahenan di oudio en bes di mensahenan di teksto si mayornan tin difikultat pa lesa. Lo bo por mester splika kon un partisipante ta skucha un mensahe di oudio."
$found29 = $d.Content.Find.Execute($find29, $true, $true, $false, $false, $false, $true, 1, $false, $replace29, 2)
if (-not $found29) { Write-Host "NOT FOUND: 29" }
